# Updated cryptos list - apply latest price/volume snapshot to sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: B=Coin, C=Link, D=Price, E=Volume(1h)
# Rows 10/11 and 42/43 and 45/46 have swapped coin identities plus new values.

# Price/Volume columns hold text that often *looks* numeric (e.g. "1.000",
# "233.72"). Force them to be stored as Text before writing so Excel does
# not silently reinterpret the strings as numbers and drop things like
# trailing zeros.
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2 - Bitcoin
$ws.Cells.Item(2, 4).Value = "30.411.18"
$ws.Cells.Item(2, 5).Value = "  +0.39%  "

# Row 3 - Ethereum
$ws.Cells.Item(3, 4).Value = "1.856.09"
$ws.Cells.Item(3, 5).Value = "  -0.19%  "

# Row 4 - TetherUSD
$ws.Cells.Item(4, 4).Value = "1.000"
$ws.Cells.Item(4, 5).Value = "  +0.05%  "

# Row 5 - BNB
$ws.Cells.Item(5, 4).Value = "233.72"
$ws.Cells.Item(5, 5).Value = "  +0.06%  "

# Row 6 - USDC
$ws.Cells.Item(6, 4).Value = "1.000"
$ws.Cells.Item(6, 5).Value = "  +0.02%  "

# Row 7 - XRP
$ws.Cells.Item(7, 4).Value = "0.4687"
$ws.Cells.Item(7, 5).Value = "  -1.49%  "

# Row 8 - Cardano
$ws.Cells.Item(8, 4).Value = "0.2741"
$ws.Cells.Item(8, 5).Value = "  -0.51%  "

# Row 9 - Dogecoin
$ws.Cells.Item(9, 4).Value = "0.06315"

# Row 10 - now Solana (was WrappedEther)
$ws.Cells.Item(10, 2).Value = "Solana"
$ws.Cells.Item(10, 3).Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Cells.Item(10, 4).Value = "17.15"
$ws.Cells.Item(10, 5).Value = "  +6.40%  "

# Row 11 - now WrappedEther (was Solana)
$ws.Cells.Item(11, 2).Value = "WrappedEther"
$ws.Cells.Item(11, 3).Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(11, 4).Value = "1.849.13"
$ws.Cells.Item(11, 5).Value = "  -0.68%  "

# Row 12 - TRON
$ws.Cells.Item(12, 4).Value = "0.07452"
$ws.Cells.Item(12, 5).Value = "  +0.46%  "

# Row 13 - Polkadot
$ws.Cells.Item(13, 4).Value = "4.941"
$ws.Cells.Item(13, 5).Value = "  -1.07%  "

# Row 14 - Litecoin
$ws.Cells.Item(14, 4).Value = "84.07"
$ws.Cells.Item(14, 5).Value = "  -1.96%  "

# Row 15 - Polygon
$ws.Cells.Item(15, 4).Value = "0.6255"
$ws.Cells.Item(15, 5).Value = "  -1.12%  "

# Row 16 - WrappedBTC
$ws.Cells.Item(16, 4).Value = "30.384.45"
$ws.Cells.Item(16, 5).Value = "  +0.35%  "

# Row 17 - Dai
$ws.Cells.Item(17, 4).Value = "0.9997"
$ws.Cells.Item(17, 5).Value = "  -0.01%  "

# Row 18 - BitcoinCash
$ws.Cells.Item(18, 4).Value = "231.03"
$ws.Cells.Item(18, 5).Value = "  +0.17%  "

# Row 19 - Avalanche
$ws.Cells.Item(19, 4).Value = "12.54"
$ws.Cells.Item(19, 5).Value = "  -1.99%  "

# Row 20 - ShibaInu
$ws.Cells.Item(20, 5).Value = "  -0.69%  "

# Row 21 - BinanceUSD
$ws.Cells.Item(21, 4).Value = "1.001"
$ws.Cells.Item(21, 5).Value = "  -0.02%  "

# Row 22 - Uniswap
$ws.Cells.Item(22, 4).Value = "4.930"
$ws.Cells.Item(22, 5).Value = "  -3.43%  "

# Row 23 - Chainlink
$ws.Cells.Item(23, 4).Value = "5.903"
$ws.Cells.Item(23, 5).Value = "  -1.94%  "

# Row 24 - Monero
$ws.Cells.Item(24, 4).Value = "167.33"
$ws.Cells.Item(24, 5).Value = "  -0.39%  "

# Row 25 - Cosmos
$ws.Cells.Item(25, 4).Value = "9.197"
$ws.Cells.Item(25, 5).Value = "  -0.83%  "

# Row 26 - EthereumClassic
$ws.Cells.Item(26, 4).Value = "17.93"
$ws.Cells.Item(26, 5).Value = "  +0.00%  "

# Row 27 - LidoDAOToken
$ws.Cells.Item(27, 4).Value = "1.878"
$ws.Cells.Item(27, 5).Value = "  +1.05%  "

# Row 28 - Stellar
$ws.Cells.Item(28, 4).Value = "0.1019"
$ws.Cells.Item(28, 5).Value = "  -0.18%  "

# Row 29 - Toncoin
$ws.Cells.Item(29, 5).Value = "  -0.59%  "

# Row 30 - InternetComputer(DFINITY)
$ws.Cells.Item(30, 4).Value = "4.084"
$ws.Cells.Item(30, 5).Value = "  -3.51%  "

# Row 31 - Filecoin
$ws.Cells.Item(31, 4).Value = "3.828"
$ws.Cells.Item(31, 5).Value = "  -2.11%  "

# Row 32 - Hedera
$ws.Cells.Item(32, 4).Value = "0.04892"
$ws.Cells.Item(32, 5).Value = "  +0.11%  "

# Row 33 - ARBITRUM
$ws.Cells.Item(33, 4).Value = "1.138"
$ws.Cells.Item(33, 5).Value = "  -0.76%  "

# Row 34 - ImmutableX
$ws.Cells.Item(34, 4).Value = "0.7035"
$ws.Cells.Item(34, 5).Value = "  -2.74%  "

# Row 35 - HuobiToken
$ws.Cells.Item(35, 4).Value = "2.713"
$ws.Cells.Item(35, 5).Value = "  +0.72%  "

# Row 36 - VeChain
$ws.Cells.Item(36, 5).Value = "  -2.37%  "

# Row 37 - MXToken
$ws.Cells.Item(37, 4).Value = "2.682"
$ws.Cells.Item(37, 5).Value = "  +1.90%  "

# Row 38 - TrustWalletToken
$ws.Cells.Item(38, 4).Value = "0.8721"
$ws.Cells.Item(38, 5).Value = "  -4.19%  "

# Row 39 - RenderToken
$ws.Cells.Item(39, 4).Value = "1.946"
$ws.Cells.Item(39, 5).Value = "  -1.80%  "

# Row 40 - Quant
$ws.Cells.Item(40, 4).Value = "105.67"
$ws.Cells.Item(40, 5).Value = "  -0.09%  "

# Row 41 - PaxDollar
$ws.Cells.Item(41, 4).Value = "1.000"
$ws.Cells.Item(41, 5).Value = "  +0.04%  "

# Row 42 - now FraxShare (was TheSandbox)
$ws.Cells.Item(42, 2).Value = "FraxShare"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(42, 4).Value = "5.527"
$ws.Cells.Item(42, 5).Value = "  -0.87%  "

# Row 43 - now TheSandbox (was FraxShare)
$ws.Cells.Item(43, 2).Value = "TheSandbox"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Cells.Item(43, 4).Value = "0.4052"
$ws.Cells.Item(43, 5).Value = "  -1.64%  "

# Row 44 - Aptos
$ws.Cells.Item(44, 4).Value = "7.178"
$ws.Cells.Item(44, 5).Value = "  +1.75%  "

# Row 45 - now Algorand (was Aave)
$ws.Cells.Item(45, 2).Value = "Algorand"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Cells.Item(45, 4).Value = "0.1227"
$ws.Cells.Item(45, 5).Value = "  +1.61%  "

# Row 46 - now Aave (was Algorand)
$ws.Cells.Item(46, 2).Value = "Aave"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Cells.Item(46, 4).Value = "61.42"
$ws.Cells.Item(46, 5).Value = "  +0.33%  "

# Row 47 - EnergySwap
$ws.Cells.Item(47, 4).Value = "8.569"
$ws.Cells.Item(47, 5).Value = "  -2.44%  "

# Row 48 - Elrond
$ws.Cells.Item(48, 4).Value = "33.44"
$ws.Cells.Item(48, 5).Value = "  +1.31%  "

# Row 49 - Cronos
$ws.Cells.Item(49, 5).Value = "  -1.03%  "

# Row 50 - NEARProtocol
$ws.Cells.Item(50, 4).Value = "1.361"
$ws.Cells.Item(50, 5).Value = "  -2.77%  "

# Row 51 - Decentraland
$ws.Cells.Item(51, 4).Value = "0.3674"
$ws.Cells.Item(51, 5).Value = "  -1.06%  "
